$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 480.40625
$ws.Range("I19").Value = 306.6154
$ws.Range("J19").Value = 599.3158
$ws.Range("K19").Value = 306.6154
$ws.Range("L19").Value = 599.3158
$ws.Range("M19").Value = -131.6154
$ws.Range("N19").Value = -949.3158

# Row 51
$ws.Range("H51").Value = 5448
$ws.Range("I51").Value = 3700.077
$ws.Range("J51").Value = 7071.0713
$ws.Range("K51").Value = 3700.077
$ws.Range("L51").Value = 7071.0713
$ws.Range("M51").Value = -3216.077
$ws.Range("N51").Value = -8039.0713

# Row 98
$ws.Range("H98").Value = 1521
$ws.Range("I98").Value = 1521
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1521
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -23

# Row 107
$ws.Range("H107").Value = 496.58066
$ws.Range("I107").Value = 211.17392
$ws.Range("J107").Value = 1317.125
$ws.Range("K107").Value = 211.17392
$ws.Range("L107").Value = 1317.125
$ws.Range("M107").Value = 1708.82608
$ws.Range("N107").Value = -5157.125

# Row 116
$ws.Range("H116").Value = 91456.586
$ws.Range("I116").Value = 135137.38
$ws.Range("J116").Value = 4095
$ws.Range("K116").Value = 135137.38
$ws.Range("L116").Value = 4095
$ws.Range("M116").Value = -131695.38
$ws.Range("N116").Value = -10979

# Row 122
$ws.Range("H122").Value = 1521
$ws.Range("I122").Value = 1521
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4563
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -2113

# Row 132
$ws.Range("H132").Value = 3489.0344
$ws.Range("I132").Value = 891.61536
$ws.Range("J132").Value = 26000
$ws.Range("K132").Value = 2674.84608
$ws.Range("L132").Value = 78000
$ws.Range("M132").Value = -144.8460800000003
$ws.Range("N132").Value = -83060

# Row 137
$ws.Range("H137").Value = 8187.125
$ws.Range("I137").Value = 9249.5
$ws.Range("K137").Value = 27748.5
$ws.Range("M137").Value = -25198.5

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1112.7142
$ws.Range("I45").Value = 886.44446
$ws.Range("J45").Value = 1520
$ws.Range("K45").Value = 886.44446
$ws.Range("L45").Value = 1520
$ws.Range("M45").Value = -509.44446
$ws.Range("N45").Value = -2274

# Row 52
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").ClearContents()
$ws.Range("N52").Value = 0

# Row 61
$ws.Range("H61").Value = 342278.6
$ws.Range("I61").Value = 264735.72
$ws.Range("K61").Value = 264735.72
$ws.Range("M61").Value = -264523.72

# Row 74
$ws.Range("H74").Value = 246226.06
$ws.Range("I74").Value = 323560
$ws.Range("J74").Value = 74986.64
$ws.Range("K74").Value = 323560
$ws.Range("L74").Value = 74986.64
$ws.Range("M74").Value = -322686
$ws.Range("N74").Value = -76734.64

# Row 77
$ws.Range("H77").Value = 246226.06
$ws.Range("I77").Value = 323560
$ws.Range("J77").Value = 74986.64
$ws.Range("K77").Value = 1617800
$ws.Range("L77").Value = 374933.2
$ws.Range("M77").Value = -1613432
$ws.Range("N77").Value = -383669.2

# Row 132
$ws.Range("H132").Value = 30321.648
$ws.Range("I132").Value = 41778.73
$ws.Range("J132").Value = 3241.2727
$ws.Range("K132").Value = 125336.19
$ws.Range("L132").Value = 9723.8181
$ws.Range("M132").Value = -122806.19
$ws.Range("N132").Value = -14783.8181

# Row 136
$ws.Range("H136").Value = 342278.6
$ws.Range("I136").Value = 264735.72
$ws.Range("K136").Value = 794207.1599999999
$ws.Range("M136").Value = -791657.1599999999

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 4777032
$ws.Range("I99").Value = 1926216.1
$ws.Range("J99").Value = 11112179
$ws.Range("K99").Value = 1926216.1
$ws.Range("L99").Value = 11112179
$ws.Range("M99").Value = -1924718.1
$ws.Range("N99").Value = -11115175

# Row 138
$ws.Range("H138").Value = 26333.334
$ws.Range("J138").Value = 26333.334
$ws.Range("L138").Value = 26333.334
$ws.Range("N138").Value = -36613.334

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2987.2666
$ws.Range("I31").Value = 1442.4166
$ws.Range("J31").Value = 9166.666999999999
$ws.Range("K31").Value = 1442.4166
$ws.Range("L31").Value = 9166.666999999999
$ws.Range("M31").Value = -1147.4166
$ws.Range("N31").Value = -9756.666999999999

# Row 34
$ws.Range("H34").Value = 2987.2666
$ws.Range("I34").Value = 1442.4166
$ws.Range("J34").Value = 9166.666999999999
$ws.Range("K34").Value = 1442.4166
$ws.Range("L34").Value = 9166.666999999999
$ws.Range("M34").Value = -1240.4166
$ws.Range("N34").Value = -9570.666999999999

# Row 94
$ws.Range("H94").Value = 5119.857
$ws.Range("I94").Value = 1307
$ws.Range("J94").Value = 8586.091
$ws.Range("K94").Value = 1307
$ws.Range("L94").Value = 8586.091
$ws.Range("M94").Value = -856
$ws.Range("N94").Value = -9488.091

# Row 110
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").ClearContents()
$ws.Range("N110").Value = 0

# Row 134
$ws.Range("H134").Value = 1990.2858
$ws.Range("I134").Value = 1031.0555
$ws.Range("K134").Value = 3093.1665
$ws.Range("M134").Value = -558.1664999999998

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 39.5
$ws.Range("I12").Value = 22.222221
$ws.Range("J12").Value = 46.260868
$ws.Range("K12").Value = 66.666663
$ws.Range("L12").Value = 138.782604
$ws.Range("M12").Value = 106.333337
$ws.Range("N12").Value = -484.782604

# Row 113
$ws.Range("H113").Value = 14706445
$ws.Range("I113").Value = 25000540
$ws.Range("J113").Value = 594.5
$ws.Range("K113").Value = 75001620
$ws.Range("L113").Value = 1783.5
$ws.Range("M113").Value = -74999450
$ws.Range("N113").Value = -6123.5

# Row 131
$ws.Range("H131").Value = 2027.0513
$ws.Range("I131").Value = 4829.3335
$ws.Range("J131").Value = 1661.5363
$ws.Range("K131").Value = 14488.0005
$ws.Range("L131").Value = 4984.6089
$ws.Range("M131").Value = -9448.000499999998
$ws.Range("N131").Value = -15064.6089

# Row 132
$ws.Range("H132").Value = 5705.4736
$ws.Range("I132").Value = 3136.7273
$ws.Range("K132").Value = 28230.5457
$ws.Range("M132").Value = -25700.5457

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 2389.7778
$ws.Range("I126").Value = 1835.2632
$ws.Range("J126").Value = 3706.75
$ws.Range("K126").Value = 5505.7896
$ws.Range("L126").Value = 11120.25
$ws.Range("M126").Value = -3035.7896
$ws.Range("N126").Value = -16060.25

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2584.6155
$ws.Range("I7").Value = 2533.3333
$ws.Range("J7").Value = 2700
$ws.Range("K7").Value = 2533.3333
$ws.Range("L7").Value = 2700
$ws.Range("M7").Value = -2421.3333
$ws.Range("N7").Value = -2924

# Row 61
$ws.Range("H61").Value = 1008.0909
$ws.Range("I61").Value = 1063
$ws.Range("J61").Value = 861.6667
$ws.Range("K61").Value = 1063
$ws.Range("L61").Value = 861.6667
$ws.Range("M61").Value = -861
$ws.Range("N61").Value = -1265.6667

# Row 113
$ws.Range("H113").Value = 1008.0909
$ws.Range("I113").Value = 1063
$ws.Range("J113").Value = 861.6667
$ws.Range("K113").Value = 1063
$ws.Range("L113").Value = 861.6667
$ws.Range("M113").Value = 1107
$ws.Range("N113").Value = -5201.6667

# Row 126
$ws.Range("H126").Value = 2584.6155
$ws.Range("I126").Value = 2533.3333
$ws.Range("J126").Value = 2700
$ws.Range("K126").Value = 7599.999899999999
$ws.Range("L126").Value = 8100
$ws.Range("M126").Value = -5129.999899999999
$ws.Range("N126").Value = -13040

# Row 136
$ws.Range("H136").Value = 4440.2104
$ws.Range("I136").Value = 2379.111
$ws.Range("K136").Value = 7137.333
$ws.Range("M136").Value = -4587.333

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 399.78262
$ws.Range("I107").Value = 346.7857
$ws.Range("J107").Value = 482.22223
$ws.Range("K107").Value = 1040.3571
$ws.Range("L107").Value = 1446.66669
$ws.Range("M107").Value = 879.6428999999998
$ws.Range("N107").Value = -5286.66669

# Row 136
$ws.Range("H136").Value = 25277572
$ws.Range("I136").Value = 37076176
$ws.Range("J136").Value = 772777.3
$ws.Range("K136").Value = 111228528
$ws.Range("L136").Value = 2318331.9
$ws.Range("M136").Value = -111225978
$ws.Range("N136").Value = -2323431.9
